# Binary Search: Aggressive Cows
# Adds two new rows (S.no. 3 & 4) to the "Binary Search 2" sheet describing
# the "Aggressive Cows" problem solved via Linear Search and Binary Search.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Pick up the same (wrap-text) cell formatting already used by row 4 so the
# new rows render consistently with the rest of the table.
$ws.Range("D4").Copy()
$ws.Range("B5:F6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the taller, wrapped row height used for the new entries.
$ws.Rows.Item(5).RowHeight = 28.8
$ws.Rows.Item(6).RowHeight = 28.8

# Row 5: Aggressive Cows solved with Linear Search.
$ws.Range("B5").Value = 3
$ws.Range("D5").Value = "Aggressive Cows - Linear Search"

# Row 6: Aggressive Cows solved with Binary Search.
$ws.Range("B6").Value = 4
$ws.Range("D6").Value = "Aggressive Cows - Binary Search"

# Shared "page" reference for both new rows.
$ws.Range("C5").Value = "Binary S2 4"
$ws.Range("C6").Value = "Binary S2 4"

# No link available for either entry.
$ws.Range("E5").Value = "na"
$ws.Range("E6").Value = "na"

# Leave the selection where the author's cursor ended up after typing the
# new rows.
$ws.Range("C7").Select()
